# "weight align fed learning" - add new experiment result sections to Sheet2
# and update a handful of existing result cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Update existing values in the "MNIST_2_20_split_label_0.300" block
# (the "3000+" text values are set further below, after the new sections,
# so that new shared-string entries are appended in the same order as the
# authored workbook)
# ---------------------------------------------------------------------------
$ws.Cells.Item(42, 3).Value = 2800
$ws.Cells.Item(43, 3).Value = 680
$ws.Cells.Item(44, 3).Value = 1052

# ---------------------------------------------------------------------------
# Fill in the iteration counts next to the "CIFAR10 20 clients iid" block
# ---------------------------------------------------------------------------
$ws.Cells.Item(50, 3).Value = 20
$ws.Cells.Item(51, 3).Value = 20
$ws.Cells.Item(52, 3).Value = 31
$ws.Cells.Item(53, 3).Value = 46
$ws.Cells.Item(54, 3).Value = 42

# ---------------------------------------------------------------------------
# Update values in the "CIFAR10 20 clients 0.3 Dir" block
# ---------------------------------------------------------------------------
$ws.Cells.Item(60, 3).Value = 144
$ws.Cells.Item(61, 3).Value = 133

# ---------------------------------------------------------------------------
# New section: CIFAR10 split label
# ---------------------------------------------------------------------------
$ws.Cells.Item(64, 1).Value = "CIFAR10 split label"

$ws.Cells.Item(65, 1).Value = "FedDyn"
$ws.Cells.Item(65, 2).Value = 85
$ws.Cells.Item(66, 1).Value = "FedDC"
$ws.Cells.Item(66, 2).Value = 85
$ws.Cells.Item(67, 1).Value = "SCAFFOLD"
$ws.Cells.Item(67, 2).Value = 85
$ws.Cells.Item(68, 1).Value = "FedProx"
$ws.Cells.Item(68, 2).Value = 85
$ws.Cells.Item(69, 1).Value = "FedAvg"
$ws.Cells.Item(69, 2).Value = 85

# ---------------------------------------------------------------------------
# New section: CIFAR100 split label
# ---------------------------------------------------------------------------
$ws.Cells.Item(72, 1).Value = "CIFAR100 split label"

$ws.Cells.Item(73, 1).Value = "FedDyn"
$ws.Cells.Item(73, 2).Value = 85
$ws.Cells.Item(74, 1).Value = "FedDC"
$ws.Cells.Item(74, 2).Value = 85
$ws.Cells.Item(75, 1).Value = "SCAFFOLD"
$ws.Cells.Item(75, 2).Value = 85
$ws.Cells.Item(76, 1).Value = "FedProx"
$ws.Cells.Item(76, 2).Value = 85
$ws.Cells.Item(77, 1).Value = "FedAvg"
$ws.Cells.Item(77, 2).Value = 85

# ---------------------------------------------------------------------------
# New section: Fashion mnist 20 dir 0.6
# ---------------------------------------------------------------------------
$ws.Cells.Item(81, 1).Value = "Fashion mnist 20 dir 0.6"

$ws.Cells.Item(82, 1).Value = "FedDyn"
$ws.Cells.Item(82, 2).Value = 89
$ws.Cells.Item(82, 3).Value = 56
$ws.Cells.Item(83, 1).Value = "FedDC"
$ws.Cells.Item(83, 2).Value = 89
$ws.Cells.Item(83, 3).Value = 422
$ws.Cells.Item(84, 1).Value = "SCAFFOLD"
$ws.Cells.Item(84, 2).Value = 89
$ws.Cells.Item(84, 3).Value = 369
$ws.Cells.Item(85, 1).Value = "FedProx"
$ws.Cells.Item(85, 2).Value = 89
$ws.Cells.Item(85, 3).Value = 320
$ws.Cells.Item(86, 1).Value = "FedAvg"
$ws.Cells.Item(86, 2).Value = 89
$ws.Cells.Item(86, 3).Value = 330

# ---------------------------------------------------------------------------
# New section: Fashion mnist 20 dir 0.3
# ---------------------------------------------------------------------------
$ws.Cells.Item(89, 1).Value = "Fashion mnist 20 dir 0.3"

$ws.Cells.Item(90, 1).Value = "FedDyn"
$ws.Cells.Item(90, 2).Value = 89
$ws.Cells.Item(91, 1).Value = "FedDC"
$ws.Cells.Item(91, 2).Value = 89
$ws.Cells.Item(92, 1).Value = "SCAFFOLD"
$ws.Cells.Item(92, 2).Value = 89
$ws.Cells.Item(93, 1).Value = "FedProx"
$ws.Cells.Item(93, 2).Value = 89
$ws.Cells.Item(94, 1).Value = "FedAvg"
$ws.Cells.Item(94, 2).Value = 89

# ---------------------------------------------------------------------------
# New section: Fashion mnist 20 iid
# ---------------------------------------------------------------------------
$ws.Cells.Item(99, 1).Value = "Fashion mnist 20 iid"

$ws.Cells.Item(100, 1).Value = "FedDyn"
$ws.Cells.Item(100, 2).Value = 89
$ws.Cells.Item(100, 3).Value = 35
$ws.Cells.Item(101, 1).Value = "FedDC"
$ws.Cells.Item(101, 2).Value = 89
$ws.Cells.Item(101, 3).Value = 100
$ws.Cells.Item(102, 1).Value = "SCAFFOLD"
$ws.Cells.Item(102, 2).Value = 89
$ws.Cells.Item(102, 3).Value = 165
$ws.Cells.Item(103, 1).Value = "FedProx"
$ws.Cells.Item(103, 2).Value = 89
$ws.Cells.Item(103, 3).Value = 151
$ws.Cells.Item(104, 1).Value = "FedAvg"
$ws.Cells.Item(104, 2).Value = 89
$ws.Cells.Item(104, 3).Value = 148

# ---------------------------------------------------------------------------
# New section: CIFAR10 20 clients 0.6 Dir (bold header, like the other
# top-level section titles in this sheet)
# ---------------------------------------------------------------------------
$ws.Cells.Item(109, 1).Value = "CIFAR10 20 clients 0.6 Dir"
$ws.Cells.Item(109, 1).Font.Bold = $true

$ws.Cells.Item(110, 1).Value = "FedDyn"
$ws.Cells.Item(110, 2).Value = 85
$ws.Cells.Item(111, 1).Value = "FedDC"
$ws.Cells.Item(111, 2).Value = 85
$ws.Cells.Item(112, 1).Value = "SCAFFOLD"
$ws.Cells.Item(112, 2).Value = 85
$ws.Cells.Item(113, 1).Value = "FedProx"
$ws.Cells.Item(113, 2).Value = 85
$ws.Cells.Item(114, 1).Value = "FedAvg"
$ws.Cells.Item(114, 2).Value = 85

# ---------------------------------------------------------------------------
# Now that all the new section headers exist, record the "3000+" iteration
# counts in the "MNIST_2_20_split_label_0.300" block (replacing "1000+")
# ---------------------------------------------------------------------------
$ws.Cells.Item(45, 3).Value = "3000+"
$ws.Cells.Item(46, 3).Value = "3000+"

# ---------------------------------------------------------------------------
# Column D is no longer used to hold long notes - narrow it back down
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.498697916666668

# ---------------------------------------------------------------------------
# Sheet view: scroll down to the new content and adjust zoom/selection
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A50").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1
$win.Zoom = 107
$ws.Range("G68").Select()
